$wb = $excel.ActiveWorkbook

# --- "Subas" sheet: the add-in keeps a running "last activity" log in A1,
#     rewritten every time the database is opened or an item is deleted.
#     Bring it to the latest recorded status. ---
$subas = $wb.Worksheets.Item("Subas")
$subas.Range("A1").Value = "Opened on 2022-01-28T12:40:46.383302400"

# --- New "Jersey" sheet/database added for the custom Sort & Search
#     feature introduced in this commit. It is appended after the
#     existing sheets and becomes the active tab. ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$jersey = $wb.Worksheets.Add($null, $lastSheet)
$jersey.Name = "Jersey"
$jersey.Activate()
